# Update the WGLC end-date text from "Ends on the 19th of July" to
# "Ends on the 29th of July" on the core WG status slide.
#
# We locate the run by its current text ("Ends on the 19") rather than by
# hard-coded slide/shape indices, so the edit is resilient to any slide
# re-ordering. Only the digits-bearing run is rewritten (via
# TextRange.Characters, matched to the exact span of the existing run) so
# the surrounding runs - the superscript "th" and " of July" - are left
# completely untouched.

$p = $ppt.ActivePresentation

$needle = "Ends on the 19"
$replacement = "Ends on the 29"

$targetSlide = $null
$targetShape = $null

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTextFrame) {
            $text = $shape.TextFrame.TextRange.Text
            if ($text.IndexOf($needle) -ge 0) {
                $targetSlide = $slide
                $targetShape = $shape
            }
        }
    }
}

if ($targetShape -ne $null) {
    $tr = $targetShape.TextFrame.TextRange
    $full = $tr.Text
    $pos = $full.IndexOf($needle) + 1
    $span = $tr.Characters($pos, $needle.Length)
    $span.Text = $replacement
}
